$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new values are plain numeric-looking strings ---
# Force text storage (NumberFormat "@") then clear the formatting
# artifact afterwards so the cell ends up with no style index,
# matching the original un-styled inline-string cells.
$numericLooking = @{
    "D5" = "412.23"
    "D6" = "130.30"
    "D7" = "0.639"
    "D9" = "0.741"
    "D11" = "43.76"
    "D12" = "0.0000226"
    "D13" = "9.39"
    "D16" = "21.28"
    "D18" = "12.68"
    "D21" = "498.76"
    "D24" = "13.61"
    "D25" = "3.42"
    "D26" = "35.09"
    "D27" = "9.22"
    "D28" = "4.80"
    "D29" = "7.65"
    "D30" = "12.21"
    "D34" = "42.27"
    "D35" = "59.37"
    "D37" = "0.0501"
    "D38" = "3.49"
    "D39" = "0.998"
    "D41" = "150.32"
    "D42" = "2.15"
    "D43" = "2.73"
    "D44" = "2.97"
    "D45" = "0.320"
    "D46" = "4.34"
    "D47" = "2.37"
    "D48" = "16.65"
    "D49" = "120.70"
    "D50" = "23.13"
}
foreach ($addr in $numericLooking.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $numericLooking.Keys) {
    $ws.Range($addr).Value = $numericLooking[$addr]
}
foreach ($addr in $numericLooking.Keys) {
    $ws.Range($addr).ClearFormats()
}

# --- Cells whose new values are already non-numeric text ---
# (URLs, coin names, multi-dot prices, percentage strings) -- these
# are never auto-coerced to numbers, so a plain assignment is safe.
$textValues = @{
    "D2" = "62.089.95"
    "E2" = "  -0.24%  "
    "D3" = "3.443.25"
    "E3" = "  +0.24%  "
    "E4" = "  -0.13%  "
    "E5" = "  +0.13%  "
    "E6" = "  +0.24%  "
    "E7" = "  +1.49%  "
    "E8" = "  +0.04%  "
    "E9" = "  -2.46%  "
    "E10" = "  +0.08%  "
    "E11" = "  +0.51%  "
    "E12" = "  +12.85%  "
    "E13" = "  +4.59%  "
    "D14" = "3.987.39"
    "E14" = "  +0.24%  "
    "E15" = "  +0.17%  "
    "E16" = "  +2.53%  "
    "D17" = "3.436.67"
    "E17" = "  +0.34%  "
    "E18" = "  +1.60%  "
    "E19" = "  +1.81%  "
    "D20" = "62.145.86"
    "E20" = "  -0.22%  "
    "E21" = "  +23.72%  "
    "E22" = "  +3.17%  "
    "E23" = "  +3.30%  "
    "E24" = "  +0.82%  "
    "E25" = "  +5.18%  "
    "E26" = "  +3.95%  "
    "E27" = "  +5.29%  "
    "E28" = "  +0.11%  "
    "E29" = "  -0.85%  "
    "E30" = "  +2.14%  "
    "E31" = "  -2.63%  "
    "E32" = "  -2.04%  "
    "E33" = "  -2.10%  "
    "E34" = "  -4.22%  "
    "E35" = "  +12.98%  "
    "E37" = "  -0.56%  "
    "B38" = "LidoDAOToken"
    "C38" = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
    "E38" = "  +2.35%  "
    "B39" = "FirstDigitalUSD"
    "C39" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "E39" = "  -0.03%  "
    "E40" = "  +4.16%  "
    "E41" = "  +6.81%  "
    "B42" = "ARBITRUM"
    "C42" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "E42" = "  +7.97%  "
    "B43" = "WEMIXToken"
    "C43" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "E43" = "  +12.86%  "
    "E44" = "  +1.75%  "
    "E45" = "  +1.21%  "
    "E46" = "  +6.09%  "
    "E47" = "  +22.01%  "
    "E48" = "  -1.06%  "
    "B49" = "BitcoinSV"
    "C49" = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
    "E49" = "  +24.40%  "
    "B50" = "EnergySwap"
    "C50" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "E50" = "  +4.00%  "
    "E51" = "  +17.72%  "
}
foreach ($addr in $textValues.Keys) {
    $ws.Range($addr).Value = $textValues[$addr]
}
